$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.778.21"
$ws.Range("E2").Value = "  +0.14%  "

$ws.Range("D3").Value = "1.642.72"
$ws.Range("E3").Value = "  -0.27%  "

$ws.Range("E4").Value = "  +0.19%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "218.44"
$ws.Range("E5").Value = "  +0.94%  "

$ws.Range("E6").Value = "  -0.75%  "

$ws.Range("E7").Value = "  +0.24%  "

$ws.Range("E8").Value = "  -0.51%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0624"
$ws.Range("E9").Value = "  -0.42%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.08"
$ws.Range("E10").Value = "  -0.71%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0846"
$ws.Range("E11").Value = "  +0.47%  "

$ws.Range("D12").Value = "1.870.15"
$ws.Range("E12").Value = "  -0.31%  "

$ws.Range("D13").Value = "1.666.95"
$ws.Range("E13").Value = "  +0.84%  "

$ws.Range("E14").Value = "  -1.10%  "

$ws.Range("E15").Value = "  -1.41%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.70"
$ws.Range("E16").Value = "  -0.98%  "

$ws.Range("D17").Value = "26.767.97"
$ws.Range("E17").Value = "  +0.07%  "

$ws.Range("D18").Value = "0.0₃0725"
$ws.Range("E18").Value = "  -2.27%  "

$ws.Range("E19").Value = "  +0.23%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "211.46"
$ws.Range("E20").Value = "  -3.01%  "

$ws.Range("E21").Value = "  -0.48%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.19"
$ws.Range("E22").Value = "  -1.08%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.32"
$ws.Range("E23").Value = "  -4.84%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.25"
$ws.Range("E24").Value = "  -2.41%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "147.52"
$ws.Range("E25").Value = "  +0.53%  "

$ws.Range("E26").Value = "  +0.23%  "

$ws.Range("E28").Value = "  -1.19%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.58"
$ws.Range("E29").Value = "  -1.02%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0502"
$ws.Range("E30").Value = "  -3.31%  "

$ws.Range("E31").Value = "  +0.85%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.35"
$ws.Range("E32").Value = "  +0.19%  "

$ws.Range("E33").Value = "  -0.64%  "

$ws.Range("D34").Value = "1.276.98"
$ws.Range("E34").Value = "  -0.18%  "

$ws.Range("E35").Value = "  -1.08%  "

$ws.Range("E36").Value = "  +0.29%  "

$ws.Range("E37").Value = "  -2.17%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.527"
$ws.Range("E38").Value = "  -1.92%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.807"
$ws.Range("E39").Value = "  -2.69%  "

$ws.Range("E40").Value = "  +0.18%  "

$ws.Range("E41").Value = "  -1.32%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.19"
$ws.Range("E42").Value = "  -2.60%  "

$ws.Range("D43").Value = "1.780.94"
$ws.Range("E43").Value = "  -0.29%  "

$ws.Range("E44").Value = "  -3.08%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "91.44"
$ws.Range("E45").Value = "  -0.65%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "60.14"
$ws.Range("E46").Value = "  +0.73%  "

$ws.Range("E47").Value = "  -1.41%  "

$ws.Range("B48").Value = "Cronos"
$ws.Range("C48").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0520"
$ws.Range("E48").Value = "  +0.88%  "

$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.55"
$ws.Range("E49").Value = "  -2.52%  "

$ws.Range("B50").Value = "Mantle"
$ws.Range("C50").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.408"
$ws.Range("E50").Value = "  -0.10%  "

$ws.Range("B51").Value = "Algorand"
$ws.Range("C51").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0960"
$ws.Range("E51").Value = "  -1.20%  "
